# Update the cryptocurrency price/volume snapshot in the worksheet.
#
# Most cells in columns D (Price) and E (Volume(1h)) hold numeric-looking
# text (e.g. "0.9479", "20.419.25") that must remain plain text rather than
# be auto-converted to numbers by Excel. Writing the value with a leading
# apostrophe (exactly as a user would type `'0.9479` into a cell) forces
# text interpretation; resetting the cell Style back to "Normal" afterwards
# avoids leaving a stray "quote prefix" style applied to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, [string]$Address, [string]$Text)
    $Sheet.Range($Address).Value = "'" + $Text
    $Sheet.Range($Address).Style = "Normal"
}

Set-TextCell $ws 'D2' '20.419.25'
Set-TextCell $ws 'E2' '  +2.42%  '
Set-TextCell $ws 'D3' '1.463.16'
Set-TextCell $ws 'D5' '0.9479'
Set-TextCell $ws 'E5' '  -5.40%  '
Set-TextCell $ws 'D6' '274.58'
Set-TextCell $ws 'E6' '  -0.55%  '
Set-TextCell $ws 'D7' '0.3650'
Set-TextCell $ws 'E7' '  -0.29%  '
Set-TextCell $ws 'D8' '0.3077'
Set-TextCell $ws 'E8' '  -0.88%  '
Set-TextCell $ws 'D9' '39.63'
Set-TextCell $ws 'D10' '1.036'
Set-TextCell $ws 'E10' '  +0.02%  '
Set-TextCell $ws 'D11' '0.06564'
Set-TextCell $ws 'E11' '  +0.74%  '
Set-TextCell $ws 'D12' '0.9990'
Set-TextCell $ws 'E12' '  -0.35%  '
Set-TextCell $ws 'D13' '17.97'
Set-TextCell $ws 'E13' '  +1.81%  '
Set-TextCell $ws 'D14' '5.409'
Set-TextCell $ws 'E14' '  -1.55%  '
Set-TextCell $ws 'D15' '6.119'
Set-TextCell $ws 'E15' '  -1.13%  '
Set-TextCell $ws 'D16' '0.00001025'
Set-TextCell $ws 'E16' '  +0.55%  '
Set-TextCell $ws 'D17' '1.460.59'
Set-TextCell $ws 'D18' '0.9645'
Set-TextCell $ws 'E18' '  -3.70%  '
Set-TextCell $ws 'D19' '0.05767'
Set-TextCell $ws 'E19' '  +1.61%  '
Set-TextCell $ws 'D20' '69.57'
Set-TextCell $ws 'E20' '  -1.98%  '
Set-TextCell $ws 'D21' '5.418'
Set-TextCell $ws 'E21' '  -3.55%  '
Set-TextCell $ws 'D22' '14.43'
Set-TextCell $ws 'E22' '  -2.06%  '
Set-TextCell $ws 'E23' '  -0.43%  '
Set-TextCell $ws 'D24' '2.236'
Set-TextCell $ws 'E24' '  +0.16%  '
Set-TextCell $ws 'D25' '20.446.77'
Set-TextCell $ws 'E25' '  +2.50%  '
Set-TextCell $ws 'D26' '141.36'
Set-TextCell $ws 'E26' '  +6.34%  '
Set-TextCell $ws 'D27' '2.081'
Set-TextCell $ws 'E27' '  -7.97%  '
Set-TextCell $ws 'D28' '17.11'
Set-TextCell $ws 'E28' '  -1.26%  '
Set-TextCell $ws 'D29' '1.613.50'
Set-TextCell $ws 'E29' '  +2.78%  '
Set-TextCell $ws 'D30' '111.90'
Set-TextCell $ws 'E30' '  +1.94%  '
Set-TextCell $ws 'D31' '3.855'
Set-TextCell $ws 'E31' '  -1.76%  '
Set-TextCell $ws 'D32' '4.875'
Set-TextCell $ws 'E32' '  -7.72%  '
Set-TextCell $ws 'B33' 'Stellar'
Set-TextCell $ws 'C33' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws 'D33' '0.07800'
Set-TextCell $ws 'E33' '  +1.30%  '
Set-TextCell $ws 'B34' 'ImmutableX'
Set-TextCell $ws 'C34' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextCell $ws 'D34' '0.7877'
Set-TextCell $ws 'E34' '  -3.55%  '
Set-TextCell $ws 'D35' '1.509'
Set-TextCell $ws 'E35' '  +2.20%  '
Set-TextCell $ws 'D36' '0.05695'
Set-TextCell $ws 'E36' '  -2.21%  '
Set-TextCell $ws 'D37' '4.660'
Set-TextCell $ws 'E37' '  -5.44%  '
Set-TextCell $ws 'D38' '1.132'
Set-TextCell $ws 'E38' '  +3.28%  '
Set-TextCell $ws 'D39' '0.02028'
Set-TextCell $ws 'E39' '  -1.69%  '
Set-TextCell $ws 'B40' 'Frax'
Set-TextCell $ws 'C40' 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextCell $ws 'D40' '0.9529'
Set-TextCell $ws 'E40' '  -4.81%  '
Set-TextCell $ws 'B41' 'Aptos'
Set-TextCell $ws 'C41' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell $ws 'D41' '10.32'
Set-TextCell $ws 'E41' '  -1.73%  '
Set-TextCell $ws 'D42' '0.1856'
Set-TextCell $ws 'E42' '  -2.13%  '
Set-TextCell $ws 'D43' '7.413'
Set-TextCell $ws 'E43' '  -10.81%  '
Set-TextCell $ws 'D44' '0.5256'
Set-TextCell $ws 'E44' '  -1.13%  '
Set-TextCell $ws 'D45' '3.485'
Set-TextCell $ws 'E45' '  -1.58%  '
Set-TextCell $ws 'D46' '11.89'
Set-TextCell $ws 'E46' '  -3.57%  '
Set-TextCell $ws 'D47' '116.99'
Set-TextCell $ws 'E47' '  +1.65%  '
Set-TextCell $ws 'D48' '0.5138'
Set-TextCell $ws 'E48' '  -0.89%  '
Set-TextCell $ws 'D49' '1.747'
Set-TextCell $ws 'E49' '  -1.45%  '
Set-TextCell $ws 'D50' '0.06420'
Set-TextCell $ws 'E50' '  +3.96%  '
Set-TextCell $ws 'D51' '0.9844'
Set-TextCell $ws 'E51' '  -1.77%  '
